$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("F1").Value = "RCI"

# Row 2 - Inhibit Scale
$ws.Range("B2").Value = 38
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 38
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 0

# Row 3 - Shift Scale
$ws.Range("B3").Value = 44
$ws.Range("C3").Value = 49
$ws.Range("D3").Value = 44
$ws.Range("E3").Value = 49
$ws.Range("F3").Value = 0

# Row 4 - Emotional Control Scale
$ws.Range("B4").Value = 39
$ws.Range("C4").Value = 21
$ws.Range("D4").Value = 39
$ws.Range("E4").Value = 21
$ws.Range("F4").Value = 0

# Row 5 - Self-Monitor Scale
$ws.Range("B5").Value = 46
$ws.Range("C5").Value = 51
$ws.Range("D5").Value = 42
$ws.Range("E5").Value = 40
$ws.Range("F5").Value = -0.5

# Row 6 - Behavioral Regulation Index
$ws.Range("B6").Value = 37
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 36
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = -0.24

# Row 7 - Initiate Scale
$ws.Range("B7").Value = 45
$ws.Range("C7").Value = 48
$ws.Range("D7").Value = 52
$ws.Range("E7").Value = 76
$ws.Range("F7").Value = 0.74

# Row 8 - Working Memory Scale
$ws.Range("B8").Value = 54
$ws.Range("C8").Value = 79
$ws.Range("D8").Value = 54
$ws.Range("E8").Value = 79
$ws.Range("F8").Value = 0

# Row 9 - Plan/Organize Scale
$ws.Range("B9").Value = 45
$ws.Range("C9").Value = 55
$ws.Range("D9").Value = 48
$ws.Range("E9").Value = 64
$ws.Range("F9").Value = 0.45

# Row 10 - Task Monitor Scale
$ws.Range("B10").Value = 42
$ws.Range("C10").Value = 45
$ws.Range("D10").Value = 42
$ws.Range("E10").Value = 45
$ws.Range("F10").Value = 0

# Row 11 - Organization of Materials Scale
$ws.Range("B11").Value = 43
$ws.Range("C11").Value = 41
$ws.Range("D11").Value = 43
$ws.Range("E11").Value = 41
$ws.Range("F11").Value = 0

# Row 12 - Metacognitive Index
$ws.Range("B12").Value = 45
$ws.Range("C12").Value = 48
$ws.Range("D12").Value = 47
$ws.Range("E12").Value = 56
$ws.Range("F12").Value = 0.48

# Row 13 - Global Executive Composite Score
$ws.Range("B13").Value = 41
$ws.Range("C13").Value = 32
$ws.Range("D13").Value = 42
$ws.Range("E13").Value = 36
$ws.Range("F13").Value = 0.31

# New column F width
$ws.Columns.Item(6).ColumnWidth = 18.67

# Update selection to match final saved state
$ws.Range("F14").Select()
